$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "29.727.01"
$ws.Range("E2").Value = "  +1.49%  "

# Row 3
$ws.Range("D3").Value = "1.995.88"
$ws.Range("E3").Value = "  +5.20%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.007"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.85%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "329.16"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.03%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.007"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.86%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4680"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.00%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3944"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.78%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "46.73"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.26%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08014"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.62%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9982"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.72%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "23.00"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.52%  "

# Row 13
$ws.Range("D13").Value = "2.026.65"
$ws.Range("E13").Value = "  +6.71%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.215"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.85%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.870"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.93%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.07135"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.17%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "88.49"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.33%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.009"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.96%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.00001004"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.38%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.41"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.60%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.004"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.67%  "

# Row 22
$ws.Range("B22").Value = "BitDAO"
$ws.Range("C22").Value = "https://coinranking.com/coin/N2IgQ9Xme+bitdao-bit"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.5194"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +10.64%  "

# Row 23
$ws.Range("B23").Value = "WrappedBTC"
$ws.Range("C23").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D23").Value = "29.782.28"
$ws.Range("E23").Value = "  +1.70%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.563"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.01%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.28"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.91%  "

# Row 26
$ws.Range("D26").Value = "2.259.70"
$ws.Range("E26").Value = "  +6.21%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.116"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.81%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "159.14"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.08%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.75"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.62%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.840"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.99%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "120.04"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.25%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.902"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.93%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09470"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.08%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9081"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.02%  "

# Row 35
$ws.Range("B35").Value = "Filecoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.276"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.49%  "

# Row 36
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.346"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.24%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.220"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.11%  "

# Row 38
$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05842"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.73%  "

# Row 39
$ws.Range("B39").Value = "TrustWalletToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.179"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.23%  "

# Row 40
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.02114"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.21%  "

# Row 41
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.927"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.69%  "

# Row 42
$ws.Range("B42").Value = "PEPE"
$ws.Range("C42").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.000003265"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +52.50%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.5770"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.03%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1822"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.07%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "9.820"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.68%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.779"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +8.69%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "11.99"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.36%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5405"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.75%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.189"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.34%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.859"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.05%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06965"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.27%  "
